$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3250
$ws.Range("I64").Value = 3241.6667
$ws.Range("J64").Value = 3300
$ws.Range("K64").Value = 3241.6667
$ws.Range("L64").Value = 3300
$ws.Range("M64").Value = -2993.6667
$ws.Range("N64").Value = -3796
$ws.Range("H67").Value = 3250
$ws.Range("I67").Value = 3241.6667
$ws.Range("J67").Value = 3300
$ws.Range("K67").Value = 3241.6667
$ws.Range("L67").Value = 3300
$ws.Range("M67").Value = -2383.6667
$ws.Range("N67").Value = -5016
$ws.Range("H137").Value = 1487.4231
$ws.Range("I137").Value = 1320.6666
$ws.Range("K137").Value = 3961.9998
$ws.Range("M137").Value = -1411.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15543.724
$ws.Range("I32").Value = 16165.819
$ws.Range("J32").Value = 4346
$ws.Range("K32").Value = 16165.819
$ws.Range("L32").Value = 4346
$ws.Range("M32").Value = -15878.819
$ws.Range("N32").Value = -4920
$ws.Range("H61").Value = 2200.1738
$ws.Range("I61").Value = 1160.4
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1160.4
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -948.4000000000001
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1011.72
$ws.Range("I74").Value = 1013.94446
$ws.Range("J74").Value = 1006
$ws.Range("K74").Value = 1013.94446
$ws.Range("L74").Value = 1006
$ws.Range("M74").Value = -139.94446
$ws.Range("N74").Value = -2754
$ws.Range("H77").Value = 1011.72
$ws.Range("I77").Value = 1013.94446
$ws.Range("J77").Value = 1006
$ws.Range("K77").Value = 5069.7223
$ws.Range("L77").Value = 5030
$ws.Range("M77").Value = -701.7223000000004
$ws.Range("N77").Value = -13766
$ws.Range("H132").Value = 4402.7295
$ws.Range("I132").Value = 6356.25
$ws.Range("J132").Value = 2104.4707
$ws.Range("K132").Value = 19068.75
$ws.Range("L132").Value = 6313.4121
$ws.Range("M132").Value = -16538.75
$ws.Range("N132").Value = -11373.4121
$ws.Range("H136").Value = 2200.1738
$ws.Range("I136").Value = 1160.4
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3481.2
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -931.2000000000003
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2966.6667
$ws.Range("I105").Value = 2966.6667
$ws.Range("K105").Value = 2966.6667
$ws.Range("M105").Value = -1219.6667
$ws.Range("H134").Value = 26105.785
$ws.Range("I134").Value = 42183.72
$ws.Range("J134").Value = 2461.7646
$ws.Range("K134").Value = 126551.16
$ws.Range("L134").Value = 7385.293799999999
$ws.Range("M134").Value = -124016.16
$ws.Range("N134").Value = -12455.2938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9525651
$ws.Range("I31").Value = 1921.2941
$ws.Range("K31").Value = 1921.2941
$ws.Range("M31").Value = -1626.2941
$ws.Range("H34").Value = 9525651
$ws.Range("I34").Value = 1921.2941
$ws.Range("K34").Value = 1921.2941
$ws.Range("M34").Value = -1719.2941
$ws.Range("H58").Value = 1604.4286
$ws.Range("I58").Value = 2640.5
$ws.Range("J58").Value = 1190
$ws.Range("K58").Value = 2640.5
$ws.Range("L58").Value = 1190
$ws.Range("M58").Value = -2437.5
$ws.Range("N58").Value = -1596
$ws.Range("H132").Value = 1945.9269
$ws.Range("I132").Value = 1781.5385
$ws.Range("J132").Value = 2230.8667
$ws.Range("K132").Value = 5344.6155
$ws.Range("L132").Value = 6692.6001
$ws.Range("M132").Value = -2814.6155
$ws.Range("N132").Value = -11752.6001
$ws.Range("H134").Value = 1264.3429
$ws.Range("I134").Value = 1223.3914
$ws.Range("J134").Value = 1342.8334
$ws.Range("K134").Value = 3670.1742
$ws.Range("L134").Value = 4028.5002
$ws.Range("M134").Value = -1135.1742
$ws.Range("N134").Value = -9098.5002
$ws.Range("H136").Value = 1604.4286
$ws.Range("I136").Value = 2640.5
$ws.Range("J136").Value = 1190
$ws.Range("K136").Value = 7921.5
$ws.Range("L136").Value = 3570
$ws.Range("M136").Value = -5371.5
$ws.Range("N136").Value = -8670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 790.5714
$ws.Range("I5").Value = 207.5
$ws.Range("J5").Value = 1568
$ws.Range("K5").Value = 622.5
$ws.Range("L5").Value = 4704
$ws.Range("M5").Value = -510.5
$ws.Range("N5").Value = -4928
$ws.Range("H135").Value = 790.5714
$ws.Range("I135").Value = 207.5
$ws.Range("J135").Value = 1568
$ws.Range("K135").Value = 1867.5
$ws.Range("L135").Value = 14112
$ws.Range("M135").Value = 667.5
$ws.Range("N135").Value = -19182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34004536
$ws.Range("I70").Value = 39235236
$ws.Range("K70").Value = 39235236
$ws.Range("M70").Value = -39234966
$ws.Range("H73").Value = 34004536
$ws.Range("I73").Value = 39235236
$ws.Range("K73").Value = 39235236
$ws.Range("M73").Value = -39234300
$ws.Range("H132").Value = 71034.62
$ws.Range("I132").Value = 134771.2
$ws.Range("J132").Value = 2745.4285
$ws.Range("K132").Value = 404313.6
$ws.Range("L132").Value = 8236.2855
$ws.Range("M132").Value = -401783.6
$ws.Range("N132").Value = -13296.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6648.974
$ws.Range("I132").Value = 8029
$ws.Range("J132").Value = 3136.182
$ws.Range("K132").Value = 24087
$ws.Range("L132").Value = 9408.545999999998
$ws.Range("M132").Value = -21557
$ws.Range("N132").Value = -14468.546
$ws.Range("H136").Value = 7567.5
$ws.Range("I136").Value = 11564.444
$ws.Range("K136").Value = 34693.33199999999
$ws.Range("M136").Value = -32143.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1296.5
$ws.Range("I132").Value = 1073.1724
$ws.Range("J132").Value = 2221.7144
$ws.Range("K132").Value = 3219.5172
$ws.Range("L132").Value = 6665.1432
$ws.Range("M132").Value = -689.5171999999998
$ws.Range("N132").Value = -11725.1432
$ws.Range("H136").Value = 8620.75
$ws.Range("I136").Value = 9128.467000000001
$ws.Range("J136").Value = 1005
$ws.Range("K136").Value = 27385.401
$ws.Range("L136").Value = 3015
$ws.Range("M136").Value = -24835.401
$ws.Range("N136").Value = -8115
